$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold text-formatted values (with
# thousands separators as dots, and percentages padded with spaces).
# For each cell being updated, force a text number-format first so Excel
# does not reinterpret / renormalize the string as a number (which would
# e.g. drop a significant trailing zero, as in "19.00" -> 19).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.403.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.454.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.75"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.453.47"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.40%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.906.99"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.274.87"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.35"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.481.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.75"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.54"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.07"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.582.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.09"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0827"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -8.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +129.44%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "430.32"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.08%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.66"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.33"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.85"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.43%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.89%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.46"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.87%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0717"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.483"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.563"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.54%  "
